# Auto-generated edit script applying market-data refresh values
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1125
$ws.Range("H108").Value = 37500
$ws.Range("J108").Value = 37500
$ws.Range("L108").Value = 37500
$ws.Range("N108").Value = -45180
$ws.Range("H129").Value = 1132.1968
$ws.Range("J129").Value = 1222.1273
$ws.Range("L129").Value = 3666.3819
$ws.Range("N129").Value = -13666.3819
$ws.Range("H137").Value = 1840.919
$ws.Range("I137").Value = 1231.65
$ws.Range("J137").Value = 2557.7058
$ws.Range("K137").Value = 3694.95
$ws.Range("L137").Value = 7673.117400000001
$ws.Range("M137").Value = -1144.95
$ws.Range("N137").Value = -12773.1174
$ws.Range("H138").Value = 3490.5493
$ws.Range("I138").Value = 2343.8696
$ws.Range("J138").Value = 4040
$ws.Range("K138").Value = 7031.6088
$ws.Range("L138").Value = 12120
$ws.Range("M138").Value = -1891.6088
$ws.Range("N138").Value = -22400
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13328.569
$ws.Range("I32").Value = 15128.298
$ws.Range("K32").Value = 15128.298
$ws.Range("M32").Value = -14841.298
$ws.Range("H74").Value = 1627.1333
$ws.Range("I74").Value = 1513.5667
$ws.Range("J74").Value = 1854.2667
$ws.Range("K74").Value = 1513.5667
$ws.Range("L74").Value = 1854.2667
$ws.Range("M74").Value = -639.5667000000001
$ws.Range("N74").Value = -3602.2667
$ws.Range("H77").Value = 1627.1333
$ws.Range("I77").Value = 1513.5667
$ws.Range("J77").Value = 1854.2667
$ws.Range("K77").Value = 7567.833500000001
$ws.Range("L77").Value = 9271.333499999999
$ws.Range("M77").Value = -3199.833500000001
$ws.Range("N77").Value = -18007.3335
$ws.Range("H110").Value = 2290.3
$ws.Range("I110").Value = 2087.875
$ws.Range("K110").Value = 2087.875
$ws.Range("M110").Value = -42.875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 119886.53
$ws.Range("I86").Value = 2357.1428
$ws.Range("K86").Value = 2357.1428
$ws.Range("M86").Value = -1234.1428
$ws.Range("H89").Value = 119886.53
$ws.Range("I89").Value = 2357.1428
$ws.Range("K89").Value = 11785.714
$ws.Range("M89").Value = -6169.714
$ws.Range("H94").Value = 46593.41
$ws.Range("I94").Value = 1105.6666
$ws.Range("J94").Value = 144067.14
$ws.Range("K94").Value = 1105.6666
$ws.Range("L94").Value = 144067.14
$ws.Range("M94").Value = -654.6666
$ws.Range("N94").Value = -144969.14
$ws.Range("H105").Value = 6213514.5
$ws.Range("I105").Value = 7520986.5
$ws.Range("J105").Value = 3021.75
$ws.Range("K105").Value = 7520986.5
$ws.Range("L105").Value = 3021.75
$ws.Range("M105").Value = -7519239.5
$ws.Range("N105").Value = -6515.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2425.6738
$ws.Range("I31").Value = 1752.9166
$ws.Range("J31").Value = 3159.5908
$ws.Range("K31").Value = 1752.9166
$ws.Range("L31").Value = 3159.5908
$ws.Range("M31").Value = -1457.9166
$ws.Range("N31").Value = -3749.5908
$ws.Range("H34").Value = 2425.6738
$ws.Range("I34").Value = 1752.9166
$ws.Range("J34").Value = 3159.5908
$ws.Range("K34").Value = 1752.9166
$ws.Range("L34").Value = 3159.5908
$ws.Range("M34").Value = -1550.9166
$ws.Range("N34").Value = -3563.5908
$ws.Range("H135").Value = 61160
$ws.Range("J135").Value = 61160
$ws.Range("L135").Value = 61160
$ws.Range("N135").Value = -71300
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 76.666664
$ws.Range("I12").Value = 76.666664
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 229.999992
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -56.99999199999999
$ws.Range("N12").ClearContents()
$ws.Range("H55").Value = 4750.5884
$ws.Range("J55").Value = 4672.5
$ws.Range("L55").Value = 14017.5
$ws.Range("N55").Value = -14371.5
$ws.Range("H58").Value = 1846.1538
$ws.Range("J58").Value = 1846.1538
$ws.Range("L58").Value = 5538.4614
$ws.Range("N58").Value = -5794.4614
$ws.Range("H68").Value = 832.73914
$ws.Range("I68").Value = 588.53125
$ws.Range("J68").Value = 1390.9286
$ws.Range("K68").Value = 1765.59375
$ws.Range("L68").Value = 4172.7858
$ws.Range("M68").Value = -954.59375
$ws.Range("N68").Value = -5794.7858
$ws.Range("H71").Value = 832.73914
$ws.Range("I71").Value = 588.53125
$ws.Range("J71").Value = 1390.9286
$ws.Range("K71").Value = 5296.78125
$ws.Range("L71").Value = 12518.3574
$ws.Range("M71").Value = -1240.78125
$ws.Range("N71").Value = -20630.3574
$ws.Range("H92").Value = 800
$ws.Range("J92").Value = 800
$ws.Range("L92").Value = 2400
$ws.Range("N92").Value = -4896
$ws.Range("H131").Value = 2235.795
$ws.Range("J131").Value = 2847.7585
$ws.Range("L131").Value = 8543.2755
$ws.Range("N131").Value = -18623.2755
$ws.Range("H132").Value = 1209.6086
$ws.Range("I132").Value = 829.8
$ws.Range("J132").Value = 1921.75
$ws.Range("K132").Value = 7468.2
$ws.Range("L132").Value = 17295.75
$ws.Range("M132").Value = -4938.2
$ws.Range("N132").Value = -22355.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3129.1428
$ws.Range("I61").Value = 1351
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 1351
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -1149
$ws.Range("N61").Value = -5904
$ws.Range("H82").Value = 1400.1428
$ws.Range("I82").Value = 832.6667
$ws.Range("J82").Value = 1825.75
$ws.Range("K82").Value = 832.6667
$ws.Range("L82").Value = 1825.75
$ws.Range("M82").Value = -471.6667
$ws.Range("N82").Value = -2547.75
$ws.Range("H85").Value = 1400.1428
$ws.Range("I85").Value = 832.6667
$ws.Range("J85").Value = 1825.75
$ws.Range("K85").Value = 832.6667
$ws.Range("L85").Value = 1825.75
$ws.Range("M85").Value = 415.3333
$ws.Range("N85").Value = -4321.75
$ws.Range("H113").Value = 3129.1428
$ws.Range("I113").Value = 1351
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 1351
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = 819
$ws.Range("N113").Value = -9840
